# Refresh market-data driven profit columns (H:N) on each class sheet.
# Mirrors the scheduled runner's nightly price/profit sync.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 1399.3462
$ws.Range("I15").Value = 1399.3462
$ws.Range("K15").Value = 4198.0386
$ws.Range("M15").Value = -4029.0386
# Row 18
$ws.Range("H18").Value = 7000.5
$ws.Range("I18").Value = 4001
$ws.Range("J18").Value = 10000
$ws.Range("K18").Value = 4001
$ws.Range("L18").Value = 10000
$ws.Range("M18").Value = -3717
$ws.Range("N18").Value = -10568
# Row 86
$ws.Range("H86").Value = 10000
$ws.Range("I86").Value = 10000
$ws.Range("K86").Value = 10000
$ws.Range("M86").Value = -8877
# Row 89
$ws.Range("H89").Value = 10000
$ws.Range("I89").Value = 10000
$ws.Range("K89").Value = 50000
$ws.Range("M89").Value = -44384
# Row 92
$ws.Range("H92").Value = 1664.3077
$ws.Range("I92").Value = 1687.3334
$ws.Range("K92").Value = 1687.3334
$ws.Range("M92").Value = -439.3334
# Row 129
$ws.Range("H129").Value = 2244.6924
$ws.Range("J129").Value = 2332
$ws.Range("L129").Value = 6996
$ws.Range("N129").Value = -16996
# Row 135
$ws.Range("H135").Value = 3694.5
$ws.Range("I135").Value = 4672.857
$ws.Range("K135").Value = 42055.713
$ws.Range("M135").Value = -39520.713
# Row 137
$ws.Range("H137").Value = 2839
$ws.Range("I137").Value = 1848.75
$ws.Range("K137").Value = 5546.25
$ws.Range("M137").Value = -2996.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3303.628
$ws.Range("I32").Value = 2677.1538
$ws.Range("K32").Value = 2677.1538
$ws.Range("M32").Value = -2390.1538
# Row 61
$ws.Range("H61").Value = 3257.75
$ws.Range("I61").Value = 3187.3333
$ws.Range("J61").Value = 3469
$ws.Range("K61").Value = 3187.3333
$ws.Range("L61").Value = 3469
$ws.Range("M61").Value = -2975.3333
$ws.Range("N61").Value = -3893
# Row 74
$ws.Range("H74").Value = 1447.5
$ws.Range("I74").Value = 797.1429000000001
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 797.1429000000001
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = 76.85709999999995
$ws.Range("N74").Value = -7748
# Row 77
$ws.Range("H77").Value = 1447.5
$ws.Range("I77").Value = 797.1429000000001
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 3985.7145
$ws.Range("L77").Value = 30000
$ws.Range("M77").Value = 382.2855
$ws.Range("N77").Value = -38736
# Row 107
$ws.Range("H107").Value = 50000
$ws.Range("J107").Value = 50000
$ws.Range("L107").Value = 50000
$ws.Range("N107").Value = -57680
# Row 136
$ws.Range("H136").Value = 3257.75
$ws.Range("I136").Value = 3187.3333
$ws.Range("J136").Value = 3469
$ws.Range("K136").Value = 9561.999899999999
$ws.Range("L136").Value = 10407
$ws.Range("M136").Value = -7011.999899999999
$ws.Range("N136").Value = -15507

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 3838.7222
$ws.Range("I134").Value = 3773.4119
$ws.Range("K134").Value = 11320.2357
$ws.Range("M134").Value = -8785.235700000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 17
$ws.Range("H17").Value = 14751.5
$ws.Range("J17").Value = 14751.5
$ws.Range("L17").Value = 14751.5
$ws.Range("N17").Value = -15099.5
# Row 86
$ws.Range("H86").Value = 35161.777
$ws.Range("I86").Value = 8600.799999999999
$ws.Range("J86").Value = 68363
$ws.Range("K86").Value = 8600.799999999999
$ws.Range("L86").Value = 68363
$ws.Range("M86").Value = -7477.799999999999
$ws.Range("N86").Value = -70609
# Row 89
$ws.Range("H89").Value = 35161.777
$ws.Range("I89").Value = 8600.799999999999
$ws.Range("J89").Value = 68363
$ws.Range("K89").Value = 43004
$ws.Range("L89").Value = 341815
$ws.Range("M89").Value = -37388
$ws.Range("N89").Value = -353047
# Row 93
$ws.Range("H93").Value = 25000
$ws.Range("I93").Value = 25000
$ws.Range("K93").Value = 25000
$ws.Range("M93").Value = -23128

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 80
$ws.Range("H80").Value = 9000.333000000001
$ws.Range("J80").Value = 8999.5
$ws.Range("L80").Value = 26998.5
$ws.Range("N80").Value = -28870.5
# Row 83
$ws.Range("H83").Value = 9000.333000000001
$ws.Range("J83").Value = 8999.5
$ws.Range("L83").Value = 80995.5
$ws.Range("N83").Value = -90355.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1840
$ws.Range("I102").Value = 1840
$ws.Range("K102").Value = 1840
$ws.Range("M102").Value = -218
# Row 105
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
# Row 126
$ws.Range("H126").Value = 2186.7693
$ws.Range("I126").Value = 2085.6
$ws.Range("J126").Value = 2250
$ws.Range("K126").Value = 6256.799999999999
$ws.Range("L126").Value = 6750
$ws.Range("M126").Value = -3786.799999999999
$ws.Range("N126").Value = -11690

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 4000
$ws.Range("J22").Value = 4000
$ws.Range("L22").Value = 4000
$ws.Range("N22").Value = -4590
# Row 27
$ws.Range("H27").Value = 4000
$ws.Range("J27").Value = 4000
$ws.Range("L27").Value = 4000
$ws.Range("N27").Value = -4214
# Row 46
$ws.Range("H46").Value = 2526
$ws.Range("I46").Value = 2526
$ws.Range("K46").Value = 2526
$ws.Range("M46").Value = -2338
# Row 55
$ws.Range("H55").Value = 596.3333
$ws.Range("I55").Value = 517.6
$ws.Range("J55").Value = 990
$ws.Range("K55").Value = 517.6
$ws.Range("L55").Value = 990
$ws.Range("M55").Value = -344.6
$ws.Range("N55").Value = -1336
# Row 132
$ws.Range("H132").Value = 5287.4
$ws.Range("I132").Value = 4225.25
$ws.Range("J132").Value = 5995.5
$ws.Range("K132").Value = 12675.75
$ws.Range("L132").Value = 17986.5
$ws.Range("M132").Value = -10145.75
$ws.Range("N132").Value = -23046.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 38906.668
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 38906.668
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 38906.668
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -39946.668
# Row 86
$ws.Range("H86").Value = 50000
$ws.Range("J86").Value = 50000
$ws.Range("L86").Value = 50000
$ws.Range("N86").Value = -52246
# Row 89
$ws.Range("H89").Value = 50000
$ws.Range("J89").Value = 50000
$ws.Range("L89").Value = 250000
$ws.Range("N89").Value = -261232
# Row 126
$ws.Range("H126").Value = 2037
$ws.Range("I126").Value = 1882.6666
$ws.Range("K126").Value = 5647.9998
$ws.Range("M126").Value = -3177.9998
# Row 132
$ws.Range("H132").Value = 2817.3
$ws.Range("I132").Value = 2484.7083
$ws.Range("K132").Value = 7454.124899999999
$ws.Range("M132").Value = -4924.124899999999
# Row 133
$ws.Range("H133").Value = 125000
$ws.Range("J133").Value = 125000
$ws.Range("L133").Value = 125000
$ws.Range("N133").Value = -135120
